$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Medstar POB North Tower -> Medstar POB South Tower (filled gaps) ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319

# --- Row 3: ownership name tweak ---
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"

# --- Row 4: address + owner + area update ---
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

# --- Row 5: area update ---
$ws.Range("L5").Value = 58717

# --- Row 6: President Madison Apartments -> Hampton House (filled gaps) ---
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580

# --- Row 7: postal code + area update ---
$ws.Range("H7").Value = 20005
$ws.Range("L7").Value = 145697

# --- Row 8: address update ---
$ws.Range("E8").Value = "1428 H ST NW"

# --- Row 10: DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis Stevens ---
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991

# "Year Built" column was previously date-formatted; the refreshed data drops
# that formatting so the values display as plain numbers.
$ws.Range("I2:I10").ClearFormats()

# Remove the now-unused "EUI Target Year" column (M) entirely.
$ws.Columns("M").Select()
$ws.Columns("M").Delete()
